$d = $word.ActiveDocument

$d.Content.Find.Execute("0.87", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.39", 2)
$d.Content.Find.Execute("0.42", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1.00", 2)
$d.Content.Find.Execute("0.33", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.75", 2)
$d.Content.Find.Execute("0.97", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.78", 2)
